$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Linked List")
$ws.Rows("3:6").Insert() | Out-Null
Write-Host "inserted"
